$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# The paragraph with "  12  " / "Distance_Range_4_point_corre_function_average"
# (split across three runs, the last one empty) becomes a single run holding
# the whole string. Find/Replace matches across run boundaries and collapses
# the match into one run automatically.
$d.Content.Find.Execute(
    "  12  Distance_Range_4_point_corre_function_average",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "  12  Distance_Range_4_point_corre_function_average", 2) | Out-Null

# --- Change 2 -------------------------------------------------------------
# Append "  Symmetry_type" to the "mfrequency ..." header line.
$d.Content.Find.Execute(
    "mfrequency  nmax  modtype   premodcoup   modcoup",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "mfrequency  nmax  modtype   premodcoup   modcoup  Symmetry_type", 2) | Out-Null

# --- Change 3 -------------------------------------------------------------
# Add the new explanatory paragraphs after that line:
#   (blank)
#   Note: for symmetry type: A1== 0, A2==1, B1==2 , B2 ==3 .
#   (blank)
#   Note: Rmax now serve as layer numbers when constructing states near initial state.
#   To compute OTOC, we also have to construct state near nearby_state. we choose 1 as
#     layer number for it because we find for 30 mode, this could grow enormous.
#     ( Rmax==4 is probably good for cyclopantane. )
#   (blank)
$newParagraphs = @(
    "",
    "Note: for symmetry type: A1== 0, A2==1, B1==2 , B2 ==3 .",
    "",
    "Note: Rmax now serve as layer numbers when constructing states near initial state.",
    "To compute OTOC, we also have to construct state near nearby_state. we choose 1 as layer number for it because we find for 30 mode, this could grow enormous. ( Rmax==4 is probably good for cyclopantane. )",
    ""
)

foreach ($text in $newParagraphs) {
    $end = $d.Paragraphs.Last.Range
    $end.Collapse(0)
    $end.InsertParagraphAfter()
    if ($text -ne "") {
        $newRange = $d.Paragraphs.Last.Range
        $newRange.InsertBefore($text)
    }
}
